$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A183").Value = "SELECT COUNT(*) AS total_july_orders FROM `"updated_table`" WHERE `"client_name`" = 'Dominos Pizza' AND DATE_TRUNC('month', `"order_created_at`") = DATE '2025-07-01' LIMIT 1;"
$ws.Range("A184").Value = "SELECT COUNT(`"id`") AS total_canceled_orders_in_july `nFROM `"updated_table`" `nWHERE `"client_name`" = 'Dominos Pizza' `nAND `"order_status`" = 'Canceled' `nAND `"final_status_at`" >= TIMESTAMP '2025-07-01 00:00:00' `nAND `"final_status_at`" < TIMESTAMP '2025-08-01 00:00:00';"
$ws.Range("A185").Value = "SELECT COUNT(`"id`") AS total_canceled_orders_in_july `nFROM `"updated_table`" `nWHERE `"client_name`" = 'Dominos Pizza' `nAND `"order_status`" = 'Canceled' `nAND `"final_status_at`" >= TIMESTAMP '2025-07-01 00:00:00' `nAND `"final_status_at`" < TIMESTAMP '2025-08-01 00:00:00';"
$ws.Range("A186").Value = "SELECT COUNT(`"id`") AS total_canceled_orders_in_july `nFROM `"updated_table`"`nWHERE `"client_name`" = 'Dominos Pizza' `nAND `"cancellation_reason`" IS NOT NULL `nAND `"final_status_at`" >= TIMESTAMP '2025-07-01 00:00:00' `nAND `"final_status_at`" < TIMESTAMP '2025-08-01 00:00:00';"
$ws.Range("A187").Value = "SELECT COUNT(`"id`") AS total_canceled_orders_in_july`nFROM `"updated_table`"`nWHERE `"client_name`" = 'Dominos Pizza'`nAND `"order_status`" = 'Canceled'`nAND `"final_status_at`" >= TIMESTAMP '2025-07-01 00:00:00'`nAND `"final_status_at`" < TIMESTAMP '2025-08-01 00:00:00'`nLIMIT 1;"
$ws.Range("A188").Value = "SELECT COUNT(`"id`") AS total_canceled_orders_in_july  `nFROM `"updated_table`" `nWHERE `"client_name`" = 'Dominos Pizza' `nAND `"order_status`" = 'Canceled'`nAND `"final_status_at`" >= TIMESTAMP '2025-07-01 00:00:00' `nAND `"final_status_at`" < TIMESTAMP '2025-08-01 00:00:00' `nLIMIT 1;"
